# Add files via upload
# Rows 21, 22, 25 and 26 become "returned" orders (Products/quantities set to
# RETURNED, total price to 0). Row 27 also becomes a "returned" order but
# keeps its original payment/date/order-type/client/address data. The
# original row 27 product data ("panadol," / "1," / "10") is moved down to a
# brand-new row 29, preceded by a brand-new "returned" row 28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every data cell in columns B:H of this sheet is stored as text, even when
# it looks numeric (e.g. "400", "0"), so force text formatting before
# assigning values to keep that convention intact.
function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# ---- Row 21: cataflam/1/400 -> RETURNED/RETURNED/0 ----
Set-TextCell 21 2 "RETURNED"
Set-TextCell 21 3 "RETURNED"
Set-TextCell 21 4 "0"

# ---- Row 22: cataflam/1/400 -> RETURNED/RETURNED/0 ----
Set-TextCell 22 2 "RETURNED"
Set-TextCell 22 3 "RETURNED"
Set-TextCell 22 4 "0"

# ---- Row 25: cataflam/10/4000 -> RETURNED/RETURNED/0 ----
Set-TextCell 25 2 "RETURNED"
Set-TextCell 25 3 "RETURNED"
Set-TextCell 25 4 "0"

# ---- Row 26: cataflam/1/400 -> RETURNED/RETURNED/0 ----
Set-TextCell 26 2 "RETURNED"
Set-TextCell 26 3 "RETURNED"
Set-TextCell 26 4 "0"

# ---- Row 27: panadol/2/200 -> RETURNED/RETURNED/0 (rest of the row, i.e.
#      payment/date/order type/client id/address, stays as it was) ----
Set-TextCell 27 2 "RETURNED"
Set-TextCell 27 3 "RETURNED"
Set-TextCell 27 4 "0"

# ---- New row 28: another RETURNED order ----
Set-TextCell 28 1 "27"
Set-TextCell 28 2 "RETURNED"
Set-TextCell 28 3 "RETURNED"
Set-TextCell 28 4 "0"
Set-TextCell 28 5 "Cash"
Set-TextCell 28 6 "2019-12-28 09:41"
Set-TextCell 28 7 "In Store"
Set-TextCell 28 8 "0"

# ---- New row 29: carries the product data that used to live in row 27 ----
Set-TextCell 29 1 "28"
Set-TextCell 29 2 "panadol,"
Set-TextCell 29 3 "1,"
Set-TextCell 29 4 "10"
Set-TextCell 29 5 "Visa"
Set-TextCell 29 6 "2019-12-28 09:41"
Set-TextCell 29 7 "In Store"
Set-TextCell 29 8 "0"

# Column A ("OrderID") is numeric throughout the sheet; restore that type for
# the two new rows (NumberFormat "@" above would otherwise make them text).
$ws.Cells.Item(28, 1).NumberFormat = "General"
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(29, 1).NumberFormat = "General"
$ws.Cells.Item(29, 1).Value = 28
